# Monthly data refresh for the AICM passenger-count report.
# A new row for "Marzo 2025" (March 2025) is inserted right above the
# existing "Febrero 2025" row (row 6), pushing every row below it down
# by one. The table/autofilter range grows by one row and the
# "Actualización: ..." footer text is bumped from Febrero to Marzo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 6 (shifts rows 6.. down to 7..).
$ws.Rows.Item(6).Insert()

# Pick up the recurring data-row formatting (style pattern s9/s10/s11)
# from the row that is now two rows below (it carries that same
# formatting cycle) and stamp it onto the freshly inserted row.
$ws.Range("B10:D10").Copy()
$ws.Range("B6:D6").PasteSpecial(-4122)

# Fill in the new March 2025 figures.
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Mar."
$ws.Range("D6").Value = 3701.671

# Grow the table (and its autofilter) by one row to keep covering the data.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B5:D92"))

# Bump the "last updated" footer caption to the new month.
$ws.Range("B93").Value = "Actualización: Marzo 2025."
